$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.25
$ws.Range("I2").Value = 2.8
$ws.Range("L2").Value = 3.6
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("Q2").Value = 2.18
$ws.Range("R2").Value = 1.69
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 12
$ws.Range("AQ2").Value = 29

# Row 4
$ws.Range("G4").Value = 1.6
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 6.25
$ws.Range("AC4").Value = 5
$ws.Range("AH4").Value = 41
$ws.Range("AJ4").Value = 7.5
$ws.Range("AK4").Value = 23
$ws.Range("AP4").Value = 21

# Row 5
$ws.Range("G5").Value = 1.67
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.3
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("S5").Value = 2.2
$ws.Range("T5").Value = 1.65
$ws.Range("U5").Value = 3.5
$ws.Range("V5").Value = 1.31
$ws.Range("W5").Value = 4
$ws.Range("X5").Value = 1.22
$ws.Range("AE5").Value = 8.5
$ws.Range("AF5").Value = 12
$ws.Range("AG5").Value = 15
$ws.Range("AI5").Value = 8
$ws.Range("AN5").Value = 12
$ws.Range("AO5").Value = 26
$ws.Range("AP5").Value = 19
$ws.Range("AQ5").Value = 67
$ws.Range("AR5").Value = 51

# Row 6
$ws.Range("G6").Value = 3.9
$ws.Range("I6").Value = 2.15
$ws.Range("J6").Value = 4.33
$ws.Range("L6").Value = 2.88
$ws.Range("Q6").Value = 1.74
$ws.Range("R6").Value = 2.1
$ws.Range("S6").Value = 2.3
$ws.Range("T6").Value = 1.6
$ws.Range("U6").Value = 3.6
$ws.Range("V6").Value = 1.3
$ws.Range("AA6").Value = 2
$ws.Range("AB6").Value = 1.75
$ws.Range("AC6").Value = 9.5
$ws.Range("AI6").Value = 7
$ws.Range("AJ6").Value = 6
$ws.Range("AK6").Value = 17
$ws.Range("AM6").Value = 401
$ws.Range("AO6").Value = 9
$ws.Range("AQ6").Value = 19
$ws.Range("AR6").Value = 19

# Row 8
$ws.Range("G8").Value = 3.3
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 2.25
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 7.5
$ws.Range("P8").Value = 2.75
$ws.Range("Q8").Value = 1.69
$ws.Range("R8").Value = 2.14
$ws.Range("S8").Value = 2.3
$ws.Range("T8").Value = 1.6
$ws.Range("U8").Value = 3.4
$ws.Range("V8").Value = 1.31
$ws.Range("W8").Value = 4.33
$ws.Range("X8").Value = 1.2
$ws.Range("Y8").Value = 1.5
$ws.Range("Z8").Value = 2.5
$ws.Range("AC8").Value = 8.5
$ws.Range("AD8").Value = 15
$ws.Range("AF8").Value = 34
$ws.Range("AG8").Value = 29
$ws.Range("AH8").Value = 41
$ws.Range("AI8").Value = 7.5
$ws.Range("AJ8").Value = 6
$ws.Range("AK8").Value = 17
$ws.Range("AL8").Value = 51
$ws.Range("AM8").Value = 800
$ws.Range("AN8").Value = 6.5
$ws.Range("AO8").Value = 10
$ws.Range("AQ8").Value = 21
$ws.Range("AR8").Value = 21
$ws.Range("AS8").Value = 34
$ws.Range("AA8").Value = 1.91
$ws.Range("AB8").Value = 1.8

# Row 10
$ws.Range("G10").Value = 2.87
$ws.Range("H10").Value = 2.77
$ws.Range("I10").Value = 2.7
$ws.Range("K10").Value = 1.88
$ws.Range("L10").Value = 3.35
$ws.Range("N10").Value = 5.2
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.32
$ws.Range("S10").Value = 2.55
$ws.Range("T10").Value = 1.45
$ws.Range("W10").Value = 4.6
$ws.Range("X10").Value = 1.16
$ws.Range("Y10").Value = 1.55
$ws.Range("Z10").Value = 2.3
$ws.Range("AA10").Value = 2.05
$ws.Range("AB10").Value = 1.7
$ws.Range("AC10").Value = 6.7
$ws.Range("AD10").Value = 13
$ws.Range("AE10").Value = 11
$ws.Range("AF10").Value = 37
$ws.Range("AG10").Value = 30
$ws.Range("AH10").Value = 50
$ws.Range("AI10").Value = 5.2
$ws.Range("AJ10").Value = 5.5
$ws.Range("AK10").Value = 17
$ws.Range("AL10").Value = 110
$ws.Range("AN10").Value = 6.4
$ws.Range("AO10").Value = 12
$ws.Range("AP10").Value = 10.5
$ws.Range("AQ10").Value = 32
$ws.Range("AR10").Value = 28
$ws.Range("AS10").Value = 45

# Row 13
$ws.Range("S13").Value = 1.9
$ws.Range("T13").Value = 1.95

# Row 14
$ws.Range("O14").Value = 1.44
$ws.Range("P14").Value = 2.75

# Row 15
$ws.Range("M15").Value = 1.07
$ws.Range("N15").Value = 9
$ws.Range("O15").Value = 1.36
$ws.Range("P15").Value = 3.2

# Row 16
$ws.Range("G16").Value = 6.5
$ws.Range("I16").Value = 1.48
$ws.Range("J16").Value = 6.5
$ws.Range("L16").Value = 2
$ws.Range("AO16").Value = 7
$ws.Range("AQ16").Value = 10

# Row 24
$ws.Range("G24").Value = 2.1
$ws.Range("H24").Value = 3.3
$ws.Range("L24").Value = 4
$ws.Range("S24").Value = 1.88
$ws.Range("T24").Value = 1.98
$ws.Range("AA24").Value = 1.7
$ws.Range("AB24").Value = 2.05
$ws.Range("AH24").Value = 26
$ws.Range("AI24").Value = 11
$ws.Range("AM24").Value = 201
$ws.Range("AN24").Value = 11

# Row 26
$ws.Range("S26").Value = 1.95
$ws.Range("T26").Value = 1.85
$ws.Range("W26").Value = 3.4
$ws.Range("X26").Value = 1.3
